# Saldo.xlsx update: remove three rows from the "Export" sheet.
#
#   - Row 2  : Conta 004480970 / ALBERTO   / 412906.51
#   - Row 7  : Conta 001882235 / LAGO      / 39557.92
#   - Row 8  : Conta 005255637 / PATRICIA  / 24053.64
#
# Deletions are performed from the bottom up (highest row number first) so
# that earlier deletions don't shift the row numbers used by later ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the LAGO / PATRICIA rows (originally rows 7-8) first.
$ws.Range("A7:A8").EntireRow.Delete()

# Then delete the ALBERTO row (row 2), which is unaffected by the deletion above.
$ws.Rows(2).Delete()
